$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the raw query-result inputs (columns B, C, D) for rows 5-13 on the
# "Query  eclipse  - Google, Bing," sheet. All the downstream formula cells
# (F:H cumulative precision, J:L cumulative recall, N:Q threshold lookups,
# the F25:H25 averages, and the cached values inside chart1/chart2) are
# driven off these inputs and recompute automatically.

$ws.Range("C5").Value = 0

$ws.Range("B6").ClearContents()

$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1

$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 1

$ws.Range("B9").ClearContents()
$ws.Range("C9").Value = 1

$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = 1

$ws.Range("B11").ClearContents()

$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0

$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1

# Move the active selection / scroll position the way the author left it.
$ws.Application.Goto($ws.Range("A10"), $false)
$ws.Range("G17").Select()
